$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Hunk 0 - row 12
$ws.Range("H12").Value = 1176.4286
$ws.Range("I12").Value = 314.55554
$ws.Range("J12").Value = 2727.8
$ws.Range("K12").Value = 314.55554
$ws.Range("L12").Value = 2727.8
$ws.Range("M12").Value = -144.55554
$ws.Range("N12").Value = -3067.8
# Hunk 1 - row 40
$ws.Range("H40").Value = 9964.52
$ws.Range("I40").Value = 7571.3
$ws.Range("J40").Value = 11560
$ws.Range("K40").Value = 7571.3
$ws.Range("L40").Value = 11560
$ws.Range("M40").Value = -7396.3
$ws.Range("N40").Value = -11910
# Hunk 2 - row 58
$ws.Range("H58").Value = 2299.8
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 2499.6667
$ws.Range("K58").Value = 6000
$ws.Range("L58").Value = 7499.000100000001
$ws.Range("M58").Value = -5850
$ws.Range("N58").Value = -7799.000100000001
# Hunk 3 - row 80
$ws.Range("H80").Value = 3429.5
$ws.Range("J80").Value = 4416.3335
$ws.Range("L80").Value = 13249.0005
$ws.Range("N80").Value = -15245.0005
# Hunk 4 - row 83
$ws.Range("H83").Value = 3429.5
$ws.Range("J83").Value = 4416.3335
$ws.Range("L83").Value = 39747.0015
$ws.Range("N83").Value = -49731.0015
# Hunk 5 - row 125
$ws.Range("H125").Value = 792.3570999999999
$ws.Range("J125").Value = 725.1111
$ws.Range("L125").Value = 6525.9999
$ws.Range("N125").Value = -11445.9999
# Hunk 6 - row 129
$ws.Range("H129").Value = 1204.75
$ws.Range("I129").Value = 931.8889
$ws.Range("J129").Value = 2023.3334
$ws.Range("K129").Value = 2795.6667
$ws.Range("L129").Value = 6070.0002
$ws.Range("M129").Value = 2204.3333
$ws.Range("N129").Value = -16070.0002
# Hunk 7 - row 131
$ws.Range("H131").Value = 6413.8237
$ws.Range("I131").Value = 6159.778
$ws.Range("J131").Value = 6699.625
$ws.Range("K131").Value = 18479.334
$ws.Range("L131").Value = 20098.875
$ws.Range("M131").Value = -13439.334
$ws.Range("N131").Value = -30178.875
# Hunk 8 - row 132
$ws.Range("H132").Value = 1701.5416
$ws.Range("I132").Value = 1701.5416
$ws.Range("K132").Value = 5104.6248
$ws.Range("M132").Value = -2574.6248
# Hunk 9 - row 135
$ws.Range("H135").Value = 782.94116
$ws.Range("I135").Value = 769.375
$ws.Range("K135").Value = 6924.375
$ws.Range("M135").Value = -4389.375
# Hunk 10 - row 137
$ws.Range("H137").Value = 3058.2563
$ws.Range("I137").Value = 2037.9231
$ws.Range("J137").Value = 5098.923
$ws.Range("K137").Value = 6113.7693
$ws.Range("L137").Value = 15296.769
$ws.Range("M137").Value = -3563.7693
$ws.Range("N137").Value = -20396.769
# Hunk 11 - row 138
$ws.Range("H138").Value = 5538.8057
$ws.Range("I138").Value = 3907.6365
$ws.Range("J138").Value = 6256.52
$ws.Range("K138").Value = 11722.9095
$ws.Range("L138").Value = 18769.56
$ws.Range("M138").Value = -6582.9095
$ws.Range("N138").Value = -29049.56

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Hunk 12 - row 61
$ws.Range("H61").Value = 6137.951
$ws.Range("I61").Value = 6165.1904
$ws.Range("J61").Value = 6109.35
$ws.Range("K61").Value = 6165.1904
$ws.Range("L61").Value = 6109.35
$ws.Range("M61").Value = -5953.1904
$ws.Range("N61").Value = -6533.35
# Hunk 13 - row 74
$ws.Range("H74").Value = 9262753
$ws.Range("I74").Value = 10104408
$ws.Range("K74").Value = 10104408
$ws.Range("M74").Value = -10103534
# Hunk 14 - row 77
$ws.Range("H77").Value = 9262753
$ws.Range("I77").Value = 10104408
$ws.Range("K77").Value = 50522040
$ws.Range("M77").Value = -50517672
# Hunk 15 - row 97
$ws.Range("H97").Value = 1888.8
$ws.Range("I97").Value = 1134.1666
$ws.Range("J97").Value = 20000
$ws.Range("K97").Value = 1134.1666
$ws.Range("L97").Value = 20000
$ws.Range("M97").Value = -638.1666
$ws.Range("N97").Value = -20992
# Hunk 16 - row 102
$ws.Range("H102").Value = 1176.3636
$ws.Range("I102").Value = 937.7778
$ws.Range("K102").Value = 937.7778
$ws.Range("M102").Value = 684.2222
# Hunk 17 - row 132
$ws.Range("H132").Value = 2287.102
$ws.Range("I132").Value = 1625.421
$ws.Range("J132").Value = 2706.1667
$ws.Range("K132").Value = 4876.263
$ws.Range("L132").Value = 8118.500100000001
$ws.Range("M132").Value = -2346.263
$ws.Range("N132").Value = -13178.5001
# Hunk 18 - row 136
$ws.Range("H136").Value = 6137.951
$ws.Range("I136").Value = 6165.1904
$ws.Range("J136").Value = 6109.35
$ws.Range("K136").Value = 18495.5712
$ws.Range("L136").Value = 18328.05
$ws.Range("M136").Value = -15945.5712
$ws.Range("N136").Value = -23428.05

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Hunk 19 - row 99
$ws.Range("H99").Value = 1959.875
$ws.Range("I99").Value = 2097
$ws.Range("K99").Value = 2097
$ws.Range("M99").Value = -599
# Hunk 20 - row 105
$ws.Range("H105").Value = 14729.6
$ws.Range("I105").Value = 20182.182
$ws.Range("J105").Value = 8065.3335
$ws.Range("K105").Value = 20182.182
$ws.Range("L105").Value = 8065.3335
$ws.Range("M105").Value = -18435.182
$ws.Range("N105").Value = -11559.3335

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Hunk 21 - row 31
$ws.Range("H31").Value = 31748.568
$ws.Range("I31").Value = 3019.35
$ws.Range("J31").Value = 65547.64999999999
$ws.Range("K31").Value = 3019.35
$ws.Range("L31").Value = 65547.64999999999
$ws.Range("M31").Value = -2724.35
$ws.Range("N31").Value = -66137.64999999999
# Hunk 22 - row 34
$ws.Range("H34").Value = 31748.568
$ws.Range("I34").Value = 3019.35
$ws.Range("J34").Value = 65547.64999999999
$ws.Range("K34").Value = 3019.35
$ws.Range("L34").Value = 65547.64999999999
$ws.Range("M34").Value = -2817.35
$ws.Range("N34").Value = -65951.64999999999
# Hunk 23 - row 58
$ws.Range("H58").Value = 6597.4375
$ws.Range("I58").Value = 3305.1428
$ws.Range("J58").Value = 9158.111000000001
$ws.Range("K58").Value = 3305.1428
$ws.Range("L58").Value = 9158.111000000001
$ws.Range("M58").Value = -3102.1428
$ws.Range("N58").Value = -9564.111000000001
# Hunk 24 - row 107
$ws.Range("H107").Value = 916.1667
$ws.Range("I107").Value = 954.06665
$ws.Range("J107").Value = 726.6667
$ws.Range("K107").Value = 954.06665
$ws.Range("L107").Value = 726.6667
$ws.Range("M107").Value = 965.93335
$ws.Range("N107").Value = -4566.6667
# Hunk 25 - row 132
$ws.Range("H132").Value = 4121.6333
$ws.Range("I132").Value = 3510.4583
$ws.Range("K132").Value = 10531.3749
$ws.Range("M132").Value = -8001.374899999999
# Hunk 26 - row 136
$ws.Range("H136").Value = 6597.4375
$ws.Range("I136").Value = 3305.1428
$ws.Range("J136").Value = 9158.111000000001
$ws.Range("K136").Value = 9915.428400000001
$ws.Range("L136").Value = 27474.333
$ws.Range("M136").Value = -7365.428400000001
$ws.Range("N136").Value = -32574.333

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Hunk 27 - row 23
$ws.Range("H23").Value = 227.375
$ws.Range("I23").Value = 43.066666
$ws.Range("J23").Value = 534.55554
$ws.Range("K23").Value = 129.199998
$ws.Range("L23").Value = 1603.66662
$ws.Range("M23").Value = 105.800002
$ws.Range("N23").Value = -2073.66662
# Hunk 28 - row 37
$ws.Range("H37").Value = 129438
$ws.Range("J37").Value = 129438
$ws.Range("L37").Value = 388314
$ws.Range("N37").Value = -388538
# Hunk 29 - row 68
$ws.Range("H68").Value = 3378.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3378.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 10135.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -11757.5
# Hunk 30 - row 70
$ws.Range("H70").Value = 10056
$ws.Range("I70").Value = 1100
$ws.Range("K70").Value = 3300
$ws.Range("M70").Value = -2985
# Hunk 31 - row 71
$ws.Range("H71").Value = 3378.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3378.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 30406.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -38518.5
# Hunk 32 - row 73
$ws.Range("H73").Value = 10056
$ws.Range("I73").Value = 1100
$ws.Range("K73").Value = 3300
$ws.Range("M73").Value = -2208
# Hunk 33 - row 131
$ws.Range("H131").Value = 9262855
$ws.Range("I131").Value = 3645.7144
$ws.Range("J131").Value = 13892460
$ws.Range("K131").Value = 10937.1432
$ws.Range("L131").Value = 41677380
$ws.Range("M131").Value = -5897.143199999999
$ws.Range("N131").Value = -41687460

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Hunk 34 - row 2
$ws.Range("H2").Value = 162.22223
$ws.Range("I2").Value = 163.75
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 163.75
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -50.75
$ws.Range("N2").Value = -376
# Hunk 35 - row 3
$ws.Range("H3").Value = 3421.4443
$ws.Range("I3").Value = 1548.8334
$ws.Range("J3").Value = 7166.6665
$ws.Range("K3").Value = 1548.8334
$ws.Range("L3").Value = 7166.6665
$ws.Range("M3").Value = -1432.8334
$ws.Range("N3").Value = -7398.6665
# Hunk 36 - row 22
$ws.Range("H22").Value = 6333.3335
$ws.Range("J22").Value = 7500
$ws.Range("L22").Value = 7500
$ws.Range("N22").Value = -8558
# Hunk 37 - row 23
$ws.Range("H23").Value = 25000
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25446
# Hunk 38 - row 25
$ws.Range("H25").Value = 1001.5
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Hunk 39 - row 22
$ws.Range("H22").Value = 8766.666999999999
$ws.Range("I22").Value = 2899
$ws.Range("J22").Value = 11700.5
$ws.Range("K22").Value = 2899
$ws.Range("L22").Value = 11700.5
$ws.Range("M22").Value = -2604
$ws.Range("N22").Value = -12290.5
# Hunk 40 - row 27
$ws.Range("H27").Value = 8766.666999999999
$ws.Range("I27").Value = 2899
$ws.Range("J27").Value = 11700.5
$ws.Range("K27").Value = 2899
$ws.Range("L27").Value = 11700.5
$ws.Range("M27").Value = -2792
$ws.Range("N27").Value = -11914.5
# Hunk 41 - row 55
$ws.Range("H55").Value = 10001000
$ws.Range("I55").Value = 12501000
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 12501000
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -12500827
$ws.Range("N55").Value = -1346
# Hunk 42 - row 68
$ws.Range("H68").Value = 4133.227
$ws.Range("I68").Value = 2759.3684
$ws.Range("K68").Value = 2759.3684
$ws.Range("M68").Value = -2010.3684
# Hunk 43 - row 71
$ws.Range("H71").Value = 4133.227
$ws.Range("I71").Value = 2759.3684
$ws.Range("K71").Value = 13796.842
$ws.Range("M71").Value = -10052.842
# Hunk 44 - row 100
$ws.Range("H100").Value = 6385.75
$ws.Range("I100").Value = 4824.7144
$ws.Range("J100").Value = 7599.8887
$ws.Range("K100").Value = 4824.7144
$ws.Range("L100").Value = 7599.8887
$ws.Range("M100").Value = -4283.7144
$ws.Range("N100").Value = -8681.8887
# Hunk 45 - row 136
$ws.Range("H136").Value = 7500.9
$ws.Range("I136").Value = 1499.8334
$ws.Range("K136").Value = 4499.5002
$ws.Range("M136").Value = -1949.5002

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Hunk 46 - row 132
$ws.Range("H132").Value = 4908.222
$ws.Range("I132").Value = 4061.16
$ws.Range("J132").Value = 6833.364
$ws.Range("K132").Value = 12183.48
$ws.Range("L132").Value = 20500.092
$ws.Range("M132").Value = -9653.48
$ws.Range("N132").Value = -25560.092
